$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Columns D that hold plain decimal-looking numbers must be forced to Text format
# so Excel does not auto-convert the string into a numeric value (which would
# drop trailing zeros / change precision), matching the source data which keeps
# these as literal text strings.
$ws.Range('D5').NumberFormat = "@"
$ws.Range('D6').NumberFormat = "@"
$ws.Range('D7').NumberFormat = "@"
$ws.Range('D8').NumberFormat = "@"
$ws.Range('D9').NumberFormat = "@"
$ws.Range('D10').NumberFormat = "@"
$ws.Range('D11').NumberFormat = "@"
$ws.Range('D12').NumberFormat = "@"
$ws.Range('D14').NumberFormat = "@"
$ws.Range('D15').NumberFormat = "@"
$ws.Range('D16').NumberFormat = "@"
$ws.Range('D18').NumberFormat = "@"
$ws.Range('D21').NumberFormat = "@"
$ws.Range('D23').NumberFormat = "@"
$ws.Range('D24').NumberFormat = "@"
$ws.Range('D25').NumberFormat = "@"
$ws.Range('D27').NumberFormat = "@"
$ws.Range('D28').NumberFormat = "@"
$ws.Range('D29').NumberFormat = "@"
$ws.Range('D30').NumberFormat = "@"
$ws.Range('D31').NumberFormat = "@"
$ws.Range('D32').NumberFormat = "@"
$ws.Range('D33').NumberFormat = "@"
$ws.Range('D34').NumberFormat = "@"
$ws.Range('D35').NumberFormat = "@"
$ws.Range('D36').NumberFormat = "@"
$ws.Range('D37').NumberFormat = "@"
$ws.Range('D38').NumberFormat = "@"
$ws.Range('D39').NumberFormat = "@"
$ws.Range('D40').NumberFormat = "@"
$ws.Range('D41').NumberFormat = "@"
$ws.Range('D42').NumberFormat = "@"
$ws.Range('D44').NumberFormat = "@"
$ws.Range('D46').NumberFormat = "@"
$ws.Range('D49').NumberFormat = "@"
$ws.Range('D50').NumberFormat = "@"
$ws.Range('D51').NumberFormat = "@"

$ws.Range('D2').Value = '27.583.60'
$ws.Range('E2').Value = '  +1.99%  '
$ws.Range('D3').Value = '1.865.82'
$ws.Range('E4').Value = '  -0.22%  '
$ws.Range('D5').Value = '312.21'
$ws.Range('E5').Value = '  +0.76%  '
$ws.Range('D6').Value = '1.013'
$ws.Range('E6').Value = '  -0.21%  '
$ws.Range('D7').Value = '0.4786'
$ws.Range('E7').Value = '  +0.55%  '
$ws.Range('D8').Value = '0.3817'
$ws.Range('E8').Value = '  +3.58%  '
$ws.Range('D9').Value = '0.07352'
$ws.Range('D10').Value = '0.9347'
$ws.Range('D11').Value = '20.82'
$ws.Range('E11').Value = '  +4.73%  '
$ws.Range('D12').Value = '0.07811'
$ws.Range('E12').Value = '  +0.35%  '
$ws.Range('D13').Value = '1.854.30'
$ws.Range('E13').Value = '  +0.64%  '
$ws.Range('D14').Value = '5.447'
$ws.Range('E14').Value = '  +0.97%  '
$ws.Range('D15').Value = '6.569'
$ws.Range('E15').Value = '  +1.34%  '
$ws.Range('D16').Value = '90.43'
$ws.Range('E16').Value = '  +1.68%  '
$ws.Range('E17').Value = '  -0.31%  '
$ws.Range('D18').Value = '0.000008840'
$ws.Range('E18').Value = '  +2.06%  '
$ws.Range('E19').Value = '  -0.24%  '
$ws.Range('D20').Value = '27.712.58'
$ws.Range('E20').Value = '  +2.42%  '
$ws.Range('D21').Value = '14.71'
$ws.Range('E21').Value = '  +1.08%  '
$ws.Range('E22').Value = '  +1.00%  '
$ws.Range('D23').Value = '10.74'
$ws.Range('E23').Value = '  +0.79%  '
$ws.Range('D24').Value = '1.939'
$ws.Range('E24').Value = '  +0.31%  '
$ws.Range('D25').Value = '156.20'
$ws.Range('E25').Value = '  +2.15%  '
$ws.Range('E26').Value = '  +0.65%  '
$ws.Range('D27').Value = '2.031'
$ws.Range('E27').Value = '  +2.15%  '
$ws.Range('D28').Value = '115.64'
$ws.Range('E28').Value = '  +0.90%  '
$ws.Range('D29').Value = '4.951'
$ws.Range('E29').Value = '  +0.49%  '
$ws.Range('D30').Value = '0.08892'
$ws.Range('E30').Value = '  +0.25%  '
$ws.Range('D31').Value = '3.330'
$ws.Range('E31').Value = '  +0.19%  '
$ws.Range('D32').Value = '1.210'
$ws.Range('E32').Value = '  +2.65%  '
$ws.Range('D33').Value = '0.7616'
$ws.Range('E33').Value = '  +3.28%  '
$ws.Range('D34').Value = '4.605'
$ws.Range('D35').Value = '2.695'
$ws.Range('E35').Value = '  +1.18%  '
$ws.Range('D36').Value = '1.134'
$ws.Range('E36').Value = '  +1.68%  '
$ws.Range('B37').Value = 'TheSandbox'
$ws.Range('C37').Value = 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand'
$ws.Range('D37').Value = '0.5714'
$ws.Range('E37').Value = '  +8.18%  '
$ws.Range('B38').Value = 'VeChain'
$ws.Range('C38').Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$ws.Range('D38').Value = '0.02036'
$ws.Range('E38').Value = '  +3.17%  '
$ws.Range('D39').Value = '0.05380'
$ws.Range('E39').Value = '  +2.26%  '
$ws.Range('D40').Value = '2.983'
$ws.Range('E40').Value = '  +0.56%  '
$ws.Range('D41').Value = '7.058'
$ws.Range('E41').Value = '  +0.28%  '
$ws.Range('D42').Value = '8.551'
$ws.Range('E42').Value = '  +3.21%  '
$ws.Range('E43').Value = '  +0.31%  '
$ws.Range('D44').Value = '0.4898'
$ws.Range('E44').Value = '  +3.27%  '
$ws.Range('E45').Value = '  +1.12%  '
$ws.Range('D46').Value = '105.49'
$ws.Range('E46').Value = '  +3.47%  '
$ws.Range('E47').Value = '  -0.26%  '
$ws.Range('E48').Value = '  +3.06%  '
$ws.Range('D49').Value = '67.50'
$ws.Range('E49').Value = '  +2.41%  '
$ws.Range('D50').Value = '0.06100'
$ws.Range('E50').Value = '  +0.60%  '
$ws.Range('D51').Value = '0.9125'
$ws.Range('E51').Value = '  +2.10%  '
